$wb = $excel.ActiveWorkbook

# ---- Sheet1: LP1912 ----
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2,1).Value = "Última actualización: 04:56:06"
$ws.Cells.Item(3,1).Value = "Total filas: 27"
$ws.Cells.Item(10,1).Value = "04:56:06"
$ws.Cells.Item(10,4).Value = 20
$ws.Cells.Item(12,1).Value = "04:56:06"
$ws.Cells.Item(12,4).Value = 26
$ws.Cells.Item(14,1).Value = "04:56:06"
$ws.Cells.Item(14,2).Value = "05:35"
$ws.Cells.Item(14,3).Value = "215B_EL PATO"
$ws.Cells.Item(14,4).Value = 39
$ws.Cells.Item(15,1).Value = "04:56:06"
$ws.Cells.Item(15,2).Value = "05:46"
$ws.Cells.Item(15,3).Value = "15_ABASTO"
$ws.Cells.Item(15,4).Value = 50
$ws.Cells.Item(16,2).Value = "05:53"
$ws.Cells.Item(16,3).Value = "10_OLMOS"
$ws.Cells.Item(16,4).Value = 71
$ws.Cells.Item(17,1).Value = "04:56:06"
$ws.Cells.Item(17,2).Value = "05:54"
$ws.Cells.Item(17,3).Value = "10_OLMOS"
$ws.Cells.Item(17,4).Value = 58
$ws.Cells.Item(18,1).Value = "04:56:06"
$ws.Cells.Item(18,2).Value = "06:04"
$ws.Cells.Item(18,3).Value = "16_SANTA ANA"
$ws.Cells.Item(18,4).Value = 68
$ws.Cells.Item(19,1).Value = "04:18:53"
$ws.Cells.Item(19,2).Value = "06:05"
$ws.Cells.Item(19,3).Value = "16_SANTA ANA"
$ws.Cells.Item(19,4).Value = 107
$ws.Cells.Item(20,2).Value = "06:11"
$ws.Cells.Item(20,3).Value = "215A_EL PATO"
$ws.Cells.Item(20,4).Value = 89
$ws.Cells.Item(21,1).Value = "04:56:06"
$ws.Cells.Item(21,2).Value = "06:12"
$ws.Cells.Item(21,3).Value = "215A_EL PATO"
$ws.Cells.Item(21,4).Value = 76
$ws.Cells.Item(22,2).Value = "06:13"
$ws.Cells.Item(22,3).Value = "225_HARAS DEL SUR"
$ws.Cells.Item(22,4).Value = 91
$ws.Cells.Item(23,1).Value = "04:56:06"
$ws.Cells.Item(23,2).Value = "06:14"
$ws.Cells.Item(23,3).Value = "225_HARAS DEL SUR"
$ws.Cells.Item(23,4).Value = 78
$ws.Cells.Item(24,1).Value = "04:42:52"
$ws.Cells.Item(24,2).Value = "06:20"
$ws.Cells.Item(24,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(24,4).Value = 98
$ws.Cells.Item(24,5).Value = "LP1912"
$ws.Cells.Item(25,1).Value = "04:56:06"
$ws.Cells.Item(25,2).Value = "06:21"
$ws.Cells.Item(25,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(25,4).Value = 85
$ws.Cells.Item(25,5).Value = "LP1912"
$ws.Cells.Item(26,1).Value = "04:42:52"
$ws.Cells.Item(26,2).Value = "06:26"
$ws.Cells.Item(26,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(26,4).Value = 104
$ws.Cells.Item(26,5).Value = "LP1912"
$ws.Cells.Item(27,1).Value = "04:56:06"
$ws.Cells.Item(27,2).Value = "06:27"
$ws.Cells.Item(27,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(27,4).Value = 91
$ws.Cells.Item(27,5).Value = "LP1912"
$ws.Cells.Item(28,1).Value = "04:42:52"
$ws.Cells.Item(28,2).Value = "06:29"
$ws.Cells.Item(28,3).Value = "86_EST CHICA-ESC AGRARIA"
$ws.Cells.Item(28,4).Value = 107
$ws.Cells.Item(28,5).Value = "LP1912"
$ws.Cells.Item(29,1).Value = "04:56:06"
$ws.Cells.Item(29,2).Value = "06:30"
$ws.Cells.Item(29,3).Value = "86_EST CHICA-ESC AGRARIA"
$ws.Cells.Item(29,4).Value = 94
$ws.Cells.Item(29,5).Value = "LP1912"
$ws.Cells.Item(30,1).Value = "04:56:06"
$ws.Cells.Item(30,2).Value = "06:31"
$ws.Cells.Item(30,3).Value = "16_SANTA ANA"
$ws.Cells.Item(30,4).Value = 95
$ws.Cells.Item(30,5).Value = "LP1912"
$ws.Cells.Item(31,1).Value = "04:56:06"
$ws.Cells.Item(31,2).Value = "06:44"
$ws.Cells.Item(31,3).Value = "225_C ROCA-H SUR"
$ws.Cells.Item(31,4).Value = 108
$ws.Cells.Item(31,5).Value = "LP1912"
$ws.Cells.Item(32,1).Value = "04:56:06"
$ws.Cells.Item(32,2).Value = "06:47"
$ws.Cells.Item(32,3).Value = "215C_EL PATO"
$ws.Cells.Item(32,4).Value = 111
$ws.Cells.Item(32,5).Value = "LP1912"

# ---- Sheet2: LP1912-215 ----
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2,1).Value = "Última actualización: 04:56:06"
$ws.Cells.Item(3,1).Value = "Total filas: 7"
$ws.Cells.Item(9,1).Value = "04:56:06"
$ws.Cells.Item(9,2).Value = "05:35"
$ws.Cells.Item(9,3).Value = "215B_EL PATO"
$ws.Cells.Item(9,4).Value = 39
$ws.Cells.Item(10,1).Value = "04:42:52"
$ws.Cells.Item(10,2).Value = "06:11"
$ws.Cells.Item(10,3).Value = "215A_EL PATO"
$ws.Cells.Item(10,4).Value = 89
$ws.Cells.Item(10,5).Value = "LP1912"
$ws.Cells.Item(11,1).Value = "04:56:06"
$ws.Cells.Item(11,2).Value = "06:12"
$ws.Cells.Item(11,3).Value = "215A_EL PATO"
$ws.Cells.Item(11,4).Value = 76
$ws.Cells.Item(11,5).Value = "LP1912"
$ws.Cells.Item(12,1).Value = "04:56:06"
$ws.Cells.Item(12,2).Value = "06:47"
$ws.Cells.Item(12,3).Value = "215C_EL PATO"
$ws.Cells.Item(12,4).Value = 111
$ws.Cells.Item(12,5).Value = "LP1912"

# ---- Sheet3: 6203-6173 ----
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2,1).Value = "Última actualización: 04:56:06"
$ws.Cells.Item(3,1).Value = "Total filas: 6"
$ws.Cells.Item(7,1).Value = "04:56:06"
$ws.Cells.Item(7,4).Value = 48
$ws.Cells.Item(9,1).Value = "04:56:06"
$ws.Cells.Item(9,2).Value = "06:09"
$ws.Cells.Item(9,3).Value = "215A_LA PLATA"
$ws.Cells.Item(9,4).Value = 73
$ws.Cells.Item(9,5).Value = "L6173"
$ws.Cells.Item(10,1).Value = "04:42:52"
$ws.Cells.Item(10,2).Value = "06:32"
$ws.Cells.Item(10,3).Value = "215C_LA PLATA"
$ws.Cells.Item(10,4).Value = 110
$ws.Cells.Item(10,5).Value = "L6203"
$ws.Cells.Item(11,1).Value = "04:56:06"
$ws.Cells.Item(11,2).Value = "06:33"
$ws.Cells.Item(11,3).Value = "215C_LA PLATA"
$ws.Cells.Item(11,4).Value = 97
$ws.Cells.Item(11,5).Value = "L6203"
